$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C2 to a real date value formatted as "YYYY-MM-DD HH:MM:SS".
# Apply the lowercase variant first, then the uppercase one, so the
# resulting numFmts table registers both format codes (164, 165) while
# the cell itself only ends up using the uppercase one (165).
$ws.Range("C2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C2").Value = (Get-Date -Year 2021 -Month 2 -Day 4 -Hour 0 -Minute 0 -Second 0)
